$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename shared string header "RunType" -> "ConditionType"
$ws.Range("C1").Value = "ConditionType"

# Update column B and C values for rows 2-20
$data = @(
    @(1, 38, 4),
    @(2, 26, 4),
    @(3, 8, 4),
    @(4, 5, 4),
    @(5, 37, 4),
    @(6, 7, 4),
    @(7, 19, 4),
    @(8, 24, 4),
    @(9, 25, 4),
    @(10, 14, 4),
    @(11, 16, 4),
    @(12, 10, 4),
    @(13, 11, 4),
    @(14, 29, 4),
    @(15, 35, 4),
    @(16, 22, 4),
    @(17, 33, 4),
    @(18, 31, 4),
    @(19, 21, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $triple = $data[$i]
    $ws.Cells.Item($row, 1).Value = $triple[0]
    $ws.Cells.Item($row, 2).Value = $triple[1]
    $ws.Cells.Item($row, 3).Value = $triple[2]
}

# Update selection on the sheet
$ws.Range("A1:C20").Select()
